$d = $word.ActiveDocument

# --- 1) "...platforms / to maximize the experience the us|er will get..." ---
# In the source doc this sentence is split as:
#   "platforms" | " to maximize the experience the us" | [[_GoBack]] | "er will get..."
# The edit retypes across the bookmark so the last two runs collapse into one
# run and the _GoBack bookmark disappears from here (it moves to the Database
# paragraph, added later). A temporary bookmark is dropped right after
# "platforms" first so the merge doesn't cascade left into that run.
$rPlatforms = $d.Content
$rPlatforms.Find.Execute("platforms", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$blockPos = $rPlatforms.End
$d.Bookmarks.Add("ZZZBlockA", $d.Range($blockPos, $blockPos)) | Out-Null

$d.Content.Find.Execute("the us" + "er", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "the user", 1) | Out-Null

$d.Bookmarks.Item("ZZZBlockA").Delete()

# --- 2) Database paragraph: split "Cards, Decks..." after "Cards, " and drop
#        the relocated _GoBack bookmark exactly there. ---
$rCards = $d.Content
$rCards.Find.Execute("Cards, Decks", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$cardsSplit = $rCards.Start + 7
$d.Bookmarks.Add("_GoBack", $d.Range($cardsSplit, $cardsSplit)) | Out-Null

# --- 3) ER diagram filename rename: "ERDiagram.pdf" -> "Picaword-ER.pdf", with
#        the result split into "For ER Diagram, please check " | "Picaword-" |
#        "ER.pdf." (and kept apart from the leading " " run before it). ---
$rBlockB = $d.Content
$rBlockB.Find.Execute("For ER Diagram, please check", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("ZZZBlockB", $d.Range($rBlockB.Start, $rBlockB.Start)) | Out-Null

$d.Content.Find.Execute("ERDiagram.pdf", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Picaword-ER.pdf", 1) | Out-Null

$d.Bookmarks.Item("ZZZBlockB").Delete()

$rSplit1 = $d.Content
$rSplit1.Find.Execute("please check Picaword", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$split1Pos = $rSplit1.Start + 13
$d.Bookmarks.Add("ZZZSplit1", $d.Range($split1Pos, $split1Pos)) | Out-Null
$d.Bookmarks.Item("ZZZSplit1").Delete()

$rSplit2 = $d.Content
$rSplit2.Find.Execute("Picaword-ER.pdf", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$split2Pos = $rSplit2.Start + 9
$d.Bookmarks.Add("ZZZSplit2", $d.Range($split2Pos, $split2Pos)) | Out-Null
$d.Bookmarks.Item("ZZZSplit2").Delete()
